# Apply the "new order" update to the SwaadSutra consolidated report.
#
# A new order (Order ID 16) came in a couple of minutes after the most
# recent existing order (Order ID 15, same customer/flat/items), so:
#   1. Duplicate the current top data row (row 2) - this both shifts every
#      existing order down by one row AND gives the brand-new row the same
#      cell typing/formatting as a real data row (instead of a freshly
#      inserted blank row).
#   2. Correct the handful of fields that differ for the new order (Order
#      ID, timestamp, requested collection date, and notes).
#   3. Update the "Daily Summary" sheet totals for 2026-01-19 to account for
#      the second order placed that day.

$wb = $excel.ActiveWorkbook

$ordersSheet = $wb.Worksheets.Item("All Orders")
$summarySheet = $wb.Worksheets.Item("Daily Summary")

# 1. Copy row 2 (order 15, the newest order) and insert the copy above
#    itself. This pushes rows 2..16 down to 3..17 and leaves row 2 as a
#    duplicate of the (now row 3) order, typed/formatted exactly like a
#    normal data row.
$ordersSheet.Rows.Item(2).Copy()
$ordersSheet.Rows.Item(2).Insert()

# 2. Correct the fields that are different for the new order 16:
#      - Order ID: 15 -> 16
#      - Date placed: 2026-01-19 05:39 -> 2026-01-19 05:41
#      - Requested collection date: 2026-01-20 -> 2026-01-21
#      - Notes: "Less Spicy" -> (cleared, no special request this time)
#    J2 is forced to Text format first so the date-looking string stays
#    literal text (matching every other Collection Date cell in the
#    sheet) instead of Excel auto-converting it to a date serial number;
#    ClearFormats afterwards drops the temporary format so the cell keeps
#    the sheet's plain default style.
$ordersSheet.Cells.Item(2, 1).Value = 16
$ordersSheet.Cells.Item(2, 2).Value = "2026-01-19 05:41"

$ordersSheet.Cells.Item(2, 10).NumberFormat = "@"
$ordersSheet.Cells.Item(2, 10).Value = "2026-01-21"
$ordersSheet.Cells.Item(2, 10).ClearFormats()

$ordersSheet.Cells.Item(2, 12).ClearContents()

# 3. Update the Daily Summary totals for 2026-01-19 (row 2) - a second order
#    was placed that day, adding its revenue (and pending amount) to the day.
$summarySheet.Cells.Item(2, 2).Value = 2
$summarySheet.Cells.Item(2, 5).Value = 210
$summarySheet.Cells.Item(2, 7).Value = 210
